# Update the "Förändrad" (Changed) date column (C) for rows 2-9
# from 2023-09-01 (45170) to 2023-09-05 (45174).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C9").Value = 45174
